$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5565830.5  # ALC!H40 (was 5858690)
$ws.Cells.Item(40, 10).Value = 10116256  # ALC!J40 (was 11127731)
$ws.Cells.Item(40, 12).Value = 10116256  # ALC!L40 (was 11127731)
$ws.Cells.Item(40, 14).Value = -10116606  # ALC!N40 (was -11128081)
$ws.Cells.Item(43, 8).Value = 5435.4546  # ALC!H43 (was 5588.9)
$ws.Cells.Item(43, 9).Value = 6989.2  # ALC!I43 (was 7749.25)
$ws.Cells.Item(43, 10).Value = 4140.6665  # ALC!J43 (was 4148.6665)
$ws.Cells.Item(43, 11).Value = 6989.2  # ALC!K43 (was 7749.25)
$ws.Cells.Item(43, 12).Value = 4140.6665  # ALC!L43 (was 4148.6665)
$ws.Cells.Item(43, 13).Value = -6920.2  # ALC!M43 (was -7680.25)
$ws.Cells.Item(43, 14).Value = -4278.6665  # ALC!N43 (was -4286.6665)
$ws.Cells.Item(80, 8).Value = 1019.55554  # ALC!H80 (was 949.2727)
$ws.Cells.Item(80, 9).Value = 1098.1666  # ALC!I80 (was 1027)
$ws.Cells.Item(80, 10).Value = 862.3333  # ALC!J80 (was 813.25)
$ws.Cells.Item(80, 11).Value = 3294.4998  # ALC!K80 (was 3081)
$ws.Cells.Item(80, 12).Value = 2586.9999  # ALC!L80 (was 2439.75)
$ws.Cells.Item(80, 13).Value = -2296.4998  # ALC!M80 (was -2083)
$ws.Cells.Item(80, 14).Value = -4582.9999  # ALC!N80 (was -4435.75)
$ws.Cells.Item(83, 8).Value = 1019.55554  # ALC!H83 (was 949.2727)
$ws.Cells.Item(83, 9).Value = 1098.1666  # ALC!I83 (was 1027)
$ws.Cells.Item(83, 10).Value = 862.3333  # ALC!J83 (was 813.25)
$ws.Cells.Item(83, 11).Value = 9883.4994  # ALC!K83 (was 9243)
$ws.Cells.Item(83, 12).Value = 7760.9997  # ALC!L83 (was 7319.25)
$ws.Cells.Item(83, 13).Value = -4891.499400000001  # ALC!M83 (was -4251)
$ws.Cells.Item(83, 14).Value = -17744.9997  # ALC!N83 (was -17303.25)
$ws.Cells.Item(96, 8).Value = 1566.4286  # ALC!H96 (was 1426.75)
$ws.Cells.Item(96, 9).Value = 1631.25  # ALC!I96 (was 1394.8)
$ws.Cells.Item(96, 11).Value = 4893.75  # ALC!K96 (was 4184.4)
$ws.Cells.Item(96, 13).Value = -3520.75  # ALC!M96 (was -2811.4)
$ws.Cells.Item(98, 8).Value = 3470.7144  # ALC!H98 (was 3259.2666)
$ws.Cells.Item(98, 9).Value = 1690.7273  # ALC!I98 (was 1574.75)
$ws.Cells.Item(98, 11).Value = 1690.7273  # ALC!K98 (was 1574.75)
$ws.Cells.Item(98, 13).Value = -192.7273  # ALC!M98 (was -76.75)
$ws.Cells.Item(122, 8).Value = 3470.7144  # ALC!H122 (was 3259.2666)
$ws.Cells.Item(122, 9).Value = 1690.7273  # ALC!I122 (was 1574.75)
$ws.Cells.Item(122, 11).Value = 5072.1819  # ALC!K122 (was 4724.25)
$ws.Cells.Item(122, 13).Value = -2622.1819  # ALC!M122 (was -2274.25)
$ws.Cells.Item(138, 8).Value = 2657.639  # ALC!H138 (was 2607.054)
$ws.Cells.Item(138, 9).Value = 1589.25  # ALC!I138 (was 1527.4615)
$ws.Cells.Item(138, 11).Value = 4767.75  # ALC!K138 (was 4582.3845)
$ws.Cells.Item(138, 13).Value = 372.25  # ALC!M138 (was 557.6154999999999)
$ws.Cells.Item(141, 8).Value = 7770.5713  # ALC!H141 (was 7674.125)
$ws.Cells.Item(141, 9).Value = 7941.5  # ALC!I141 (was 7806.857)
$ws.Cells.Item(141, 11).Value = 23824.5  # ALC!K141 (was 23420.571)
$ws.Cells.Item(141, 13).Value = -18644.5  # ALC!M141 (was -18240.571)

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4490.3794  # ARM!H32 (was 4094.0938)
$ws.Cells.Item(32, 9).Value = 2523.25  # ARM!I32 (was 2272.1482)
$ws.Cells.Item(32, 11).Value = 2523.25  # ARM!K32 (was 2272.1482)
$ws.Cells.Item(32, 13).Value = -2236.25  # ARM!M32 (was -1985.1482)
$ws.Cells.Item(61, 8).Value = 55557390  # ARM!H61 (was 58825380)
$ws.Cells.Item(61, 9).Value = 62501840  # ARM!I61 (was 66668524)
$ws.Cells.Item(61, 11).Value = 62501840  # ARM!K61 (was 66668524)
$ws.Cells.Item(61, 13).Value = -62501628  # ARM!M61 (was -66668312)
$ws.Cells.Item(74, 8).Value = 35719424  # ARM!H74 (was 38466988)
$ws.Cells.Item(74, 9).Value = 37042324  # ARM!I74 (was 38466988)
$ws.Cells.Item(74, 10).Value = 1100  # ARM!J74 (was 0)
$ws.Cells.Item(74, 11).Value = 37042324  # ARM!K74 (was 38466988)
$ws.Cells.Item(74, 12).Value = 1100  # ARM!L74 (was 0)
$ws.Cells.Item(74, 13).Value = -37041450  # ARM!M74 (was -38466114)
$ws.Cells.Item(74, 14).Value = -2848  # ARM!N74 (was None)
$ws.Cells.Item(77, 8).Value = 35719424  # ARM!H77 (was 38466988)
$ws.Cells.Item(77, 9).Value = 37042324  # ARM!I77 (was 38466988)
$ws.Cells.Item(77, 10).Value = 1100  # ARM!J77 (was 0)
$ws.Cells.Item(77, 11).Value = 185211620  # ARM!K77 (was 192334940)
$ws.Cells.Item(77, 12).Value = 5500  # ARM!L77 (was 0)
$ws.Cells.Item(77, 13).Value = -185207252  # ARM!M77 (was -192330572)
$ws.Cells.Item(77, 14).Value = -14236  # ARM!N77 (was None)
$ws.Cells.Item(102, 8).Value = 11112086  # ARM!H102 (was 10000954)
$ws.Cells.Item(102, 9).Value = 11112086  # ARM!I102 (was 10000954)
$ws.Cells.Item(102, 11).Value = 11112086  # ARM!K102 (was 10000954)
$ws.Cells.Item(102, 13).Value = -11110464  # ARM!M102 (was -9999332)
$ws.Cells.Item(110, 8).Value = 78203.305  # ARM!H110 (was 84649.5)
$ws.Cells.Item(110, 9).Value = 91976.63  # ARM!I110 (was 101089.4)
$ws.Cells.Item(110, 11).Value = 91976.63  # ARM!K110 (was 101089.4)
$ws.Cells.Item(110, 13).Value = -89931.63  # ARM!M110 (was -99044.4)
$ws.Cells.Item(136, 8).Value = 55557390  # ARM!H136 (was 58825380)
$ws.Cells.Item(136, 9).Value = 62501840  # ARM!I136 (was 66668524)
$ws.Cells.Item(136, 11).Value = 187505520  # ARM!K136 (was 200005572)
$ws.Cells.Item(136, 13).Value = -187502970  # ARM!M136 (was -200003022)

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2605.0527  # BSM!H20 (was 2805.2222)
$ws.Cells.Item(20, 9).Value = 2814  # BSM!I20 (was 3059.7)
$ws.Cells.Item(20, 10).Value = 2317.75  # BSM!J20 (was 2487.125)
$ws.Cells.Item(20, 11).Value = 2814  # BSM!K20 (was 3059.7)
$ws.Cells.Item(20, 12).Value = 2317.75  # BSM!L20 (was 2487.125)
$ws.Cells.Item(20, 13).Value = -2567  # BSM!M20 (was -2812.7)
$ws.Cells.Item(20, 14).Value = -2811.75  # BSM!N20 (was -2981.125)
$ws.Cells.Item(64, 8).Value = 513.25  # BSM!H64 (was 528.6667)
$ws.Cells.Item(64, 10).Value = 468  # BSM!J64 (was 469)
$ws.Cells.Item(64, 12).Value = 468  # BSM!L64 (was 469)
$ws.Cells.Item(64, 14).Value = -918  # BSM!N64 (was -919)
$ws.Cells.Item(67, 8).Value = 513.25  # BSM!H67 (was 528.6667)
$ws.Cells.Item(67, 10).Value = 468  # BSM!J67 (was 469)
$ws.Cells.Item(67, 12).Value = 468  # BSM!L67 (was 469)
$ws.Cells.Item(67, 14).Value = -2028  # BSM!N67 (was -2029)
$ws.Cells.Item(86, 8).Value = 1842.2106  # BSM!H86 (was 1890.4595)
$ws.Cells.Item(86, 9).Value = 1795.3667  # BSM!I86 (was 1808.4667)
$ws.Cells.Item(86, 10).Value = 2017.875  # BSM!J86 (was 2241.8572)
$ws.Cells.Item(86, 11).Value = 1795.3667  # BSM!K86 (was 1808.4667)
$ws.Cells.Item(86, 12).Value = 2017.875  # BSM!L86 (was 2241.8572)
$ws.Cells.Item(86, 13).Value = -672.3667  # BSM!M86 (was -685.4667)
$ws.Cells.Item(86, 14).Value = -4263.875  # BSM!N86 (was -4487.8572)
$ws.Cells.Item(89, 8).Value = 1842.2106  # BSM!H89 (was 1890.4595)
$ws.Cells.Item(89, 9).Value = 1795.3667  # BSM!I89 (was 1808.4667)
$ws.Cells.Item(89, 10).Value = 2017.875  # BSM!J89 (was 2241.8572)
$ws.Cells.Item(89, 11).Value = 8976.8335  # BSM!K89 (was 9042.3335)
$ws.Cells.Item(89, 12).Value = 10089.375  # BSM!L89 (was 11209.286)
$ws.Cells.Item(89, 13).Value = -3360.833500000001  # BSM!M89 (was -3426.333500000001)
$ws.Cells.Item(89, 14).Value = -21321.375  # BSM!N89 (was -22441.286)
$ws.Cells.Item(94, 8).Value = 917.4  # BSM!H94 (was 753.8889)
$ws.Cells.Item(94, 9).Value = 917.4  # BSM!I94 (was 831.1667)
$ws.Cells.Item(94, 10).Value = 0  # BSM!J94 (was 599.3333)
$ws.Cells.Item(94, 11).Value = 917.4  # BSM!K94 (was 831.1667)
$ws.Cells.Item(94, 12).Value = 0  # BSM!L94 (was 599.3333)
$ws.Cells.Item(94, 13).Value = -466.4  # BSM!M94 (was -380.1667)
$ws.Cells.Item(94, 14).Value = $null  # BSM!N94 clear (was -1501.3333)
$ws.Cells.Item(105, 8).Value = 2440.3157  # BSM!H105 (was 2200.9412)
$ws.Cells.Item(105, 9).Value = 1883.2858  # BSM!I105 (was 1451.3334)
$ws.Cells.Item(105, 11).Value = 1883.2858  # BSM!K105 (was 1451.3334)
$ws.Cells.Item(105, 13).Value = -136.2858000000001  # BSM!M105 (was 295.6666)
$ws.Cells.Item(107, 8).Value = 128895.125  # BSM!H107 (was 45644.566)
$ws.Cells.Item(107, 9).Value = 4451.7144  # BSM!I107 (was 2326.3)
$ws.Cells.Item(107, 10).Value = 999999  # BSM!J107 (was 334433)
$ws.Cells.Item(107, 11).Value = 4451.7144  # BSM!K107 (was 2326.3)
$ws.Cells.Item(107, 12).Value = 999999  # BSM!L107 (was 334433)
$ws.Cells.Item(107, 13).Value = -2531.7144  # BSM!M107 (was -406.3000000000002)
$ws.Cells.Item(107, 14).Value = -1003839  # BSM!N107 (was -338273)

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7165.3057  # CRP!H31 (was 7701.9395)
$ws.Cells.Item(31, 9).Value = 4869.25  # CRP!I31 (was 5212.0454)
$ws.Cells.Item(31, 10).Value = 11757.417  # CRP!J31 (was 12681.728)
$ws.Cells.Item(31, 11).Value = 4869.25  # CRP!K31 (was 5212.0454)
$ws.Cells.Item(31, 12).Value = 11757.417  # CRP!L31 (was 12681.728)
$ws.Cells.Item(31, 13).Value = -4574.25  # CRP!M31 (was -4917.0454)
$ws.Cells.Item(31, 14).Value = -12347.417  # CRP!N31 (was -13271.728)
$ws.Cells.Item(33, 8).Value = 5557  # CRP!H33 (was 5946.25)
$ws.Cells.Item(33, 9).Value = 5557  # CRP!I33 (was 5946.25)
$ws.Cells.Item(33, 11).Value = 5557  # CRP!K33 (was 5946.25)
$ws.Cells.Item(33, 13).Value = -5178  # CRP!M33 (was -5567.25)
$ws.Cells.Item(34, 8).Value = 7165.3057  # CRP!H34 (was 7701.9395)
$ws.Cells.Item(34, 9).Value = 4869.25  # CRP!I34 (was 5212.0454)
$ws.Cells.Item(34, 10).Value = 11757.417  # CRP!J34 (was 12681.728)
$ws.Cells.Item(34, 11).Value = 4869.25  # CRP!K34 (was 5212.0454)
$ws.Cells.Item(34, 12).Value = 11757.417  # CRP!L34 (was 12681.728)
$ws.Cells.Item(34, 13).Value = -4667.25  # CRP!M34 (was -5010.0454)
$ws.Cells.Item(34, 14).Value = -12161.417  # CRP!N34 (was -13085.728)
$ws.Cells.Item(58, 8).Value = 20839126  # CRP!H58 (was 20839130)
$ws.Cells.Item(58, 9).Value = 29418954  # CRP!I58 (was 29418960)
$ws.Cells.Item(58, 10).Value = 2401  # CRP!J58 (was 2401.2856)
$ws.Cells.Item(58, 11).Value = 29418954  # CRP!K58 (was 29418960)
$ws.Cells.Item(58, 12).Value = 2401  # CRP!L58 (was 2401.2856)
$ws.Cells.Item(58, 13).Value = -29418751  # CRP!M58 (was -29418757)
$ws.Cells.Item(58, 14).Value = -2807  # CRP!N58 (was -2807.2856)
$ws.Cells.Item(107, 8).Value = 72409.93  # CRP!H107 (was 72385.71)
$ws.Cells.Item(107, 9).Value = 682  # CRP!I107 (was 615.625)
$ws.Cells.Item(107, 10).Value = 144137.86  # CRP!J107 (was 168079.17)
$ws.Cells.Item(107, 11).Value = 682  # CRP!K107 (was 615.625)
$ws.Cells.Item(107, 12).Value = 144137.86  # CRP!L107 (was 168079.17)
$ws.Cells.Item(107, 13).Value = 1238  # CRP!M107 (was 1304.375)
$ws.Cells.Item(107, 14).Value = -147977.86  # CRP!N107 (was -171919.17)
$ws.Cells.Item(136, 8).Value = 20839126  # CRP!H136 (was 20839130)
$ws.Cells.Item(136, 9).Value = 29418954  # CRP!I136 (was 29418960)
$ws.Cells.Item(136, 10).Value = 2401  # CRP!J136 (was 2401.2856)
$ws.Cells.Item(136, 11).Value = 88256862  # CRP!K136 (was 88256880)
$ws.Cells.Item(136, 12).Value = 7203  # CRP!L136 (was 7203.8568)
$ws.Cells.Item(136, 13).Value = -88254312  # CRP!M136 (was -88254330)
$ws.Cells.Item(136, 14).Value = -12303  # CRP!N136 (was -12303.8568)

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 2929.2856  # CUL!H80 (was 3278.8)
$ws.Cells.Item(80, 9).Value = 2701.6667  # CUL!I80 (was 3497)
$ws.Cells.Item(80, 10).Value = 3100  # CUL!J80 (was 3133.3333)
$ws.Cells.Item(80, 11).Value = 8105.000100000001  # CUL!K80 (was 10491)
$ws.Cells.Item(80, 12).Value = 9300  # CUL!L80 (was 9399.999899999999)
$ws.Cells.Item(80, 13).Value = -7169.000100000001  # CUL!M80 (was -9555)
$ws.Cells.Item(80, 14).Value = -11172  # CUL!N80 (was -11271.9999)
$ws.Cells.Item(83, 8).Value = 2929.2856  # CUL!H83 (was 3278.8)
$ws.Cells.Item(83, 9).Value = 2701.6667  # CUL!I83 (was 3497)
$ws.Cells.Item(83, 10).Value = 3100  # CUL!J83 (was 3133.3333)
$ws.Cells.Item(83, 11).Value = 24315.0003  # CUL!K83 (was 31473)
$ws.Cells.Item(83, 12).Value = 27900  # CUL!L83 (was 28199.9997)
$ws.Cells.Item(83, 13).Value = -19635.0003  # CUL!M83 (was -26793)
$ws.Cells.Item(83, 14).Value = -37260  # CUL!N83 (was -37559.9997)
$ws.Cells.Item(118, 8).Value = 1639.5  # CUL!H118 (was 2500)
$ws.Cells.Item(118, 9).Value = 1352.6666  # CUL!I118 (was 0)
$ws.Cells.Item(118, 11).Value = 4057.9998  # CUL!K118 (was 0)
$ws.Cells.Item(118, 13).Value = -2814.9998  # CUL!M118 (was None)
$ws.Cells.Item(131, 8).Value = 1659.35  # CUL!H131 (was 1641.5264)
$ws.Cells.Item(131, 10).Value = 2221.7778  # CUL!J131 (was 2249.75)
$ws.Cells.Item(131, 12).Value = 6665.3334  # CUL!L131 (was 6749.25)
$ws.Cells.Item(131, 14).Value = -16745.3334  # CUL!N131 (was -16829.25)
$ws.Cells.Item(132, 8).Value = 6666.3335  # CUL!H132 (was 20249.75)
$ws.Cells.Item(132, 9).Value = 4999  # CUL!I132 (was 32999.5)
$ws.Cells.Item(132, 11).Value = 44991  # CUL!K132 (was 296995.5)
$ws.Cells.Item(132, 13).Value = -42461  # CUL!M132 (was -294465.5)

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 9999999  # GSM!H14 (was 5834066)
$ws.Cells.Item(14, 9).Value = 9999999  # GSM!I14 (was 5834066)
$ws.Cells.Item(14, 11).Value = 9999999  # GSM!K14 (was 5834066)
$ws.Cells.Item(14, 13).Value = -9999831  # GSM!M14 (was -5833898)
$ws.Cells.Item(15, 8).Value = 59999  # GSM!H15 (was 52989.5)
$ws.Cells.Item(15, 10).Value = 59999  # GSM!J15 (was 52989.5)
$ws.Cells.Item(15, 12).Value = 59999  # GSM!L15 (was 52989.5)
$ws.Cells.Item(15, 14).Value = -60575  # GSM!N15 (was -53565.5)
$ws.Cells.Item(29, 8).Value = 16247.25  # GSM!H29 (was 14909.5)
$ws.Cells.Item(29, 9).Value = 4999.6665  # GSM!I29 (was 7893.4)
$ws.Cells.Item(29, 11).Value = 4999.6665  # GSM!K29 (was 7893.4)
$ws.Cells.Item(29, 13).Value = -4709.6665  # GSM!M29 (was -7603.4)
$ws.Cells.Item(39, 8).Value = 44799.6  # GSM!H39 (was 45666.668)
$ws.Cells.Item(39, 10).Value = 44799.6  # GSM!J39 (was 45666.668)
$ws.Cells.Item(39, 12).Value = 44799.6  # GSM!L39 (was 45666.668)
$ws.Cells.Item(39, 14).Value = -45863.6  # GSM!N39 (was -46730.668)
$ws.Cells.Item(80, 10).Value = 2000  # GSM!J80 (was 0)
$ws.Cells.Item(80, 12).Value = 2000  # GSM!L80 (was 0)
$ws.Cells.Item(80, 14).Value = -3996  # GSM!N80 (was None)
$ws.Cells.Item(81, 8).Value = 59999  # GSM!H81 (was 52989.5)
$ws.Cells.Item(81, 10).Value = 59999  # GSM!J81 (was 52989.5)
$ws.Cells.Item(81, 12).Value = 59999  # GSM!L81 (was 52989.5)
$ws.Cells.Item(81, 14).Value = -61995  # GSM!N81 (was -54985.5)
$ws.Cells.Item(83, 10).Value = 2000  # GSM!J83 (was 0)
$ws.Cells.Item(83, 12).Value = 10000  # GSM!L83 (was 0)
$ws.Cells.Item(83, 14).Value = -19984  # GSM!N83 (was None)
$ws.Cells.Item(84, 8).Value = 59999  # GSM!H84 (was 52989.5)
$ws.Cells.Item(84, 10).Value = 59999  # GSM!J84 (was 52989.5)
$ws.Cells.Item(84, 12).Value = 179997  # GSM!L84 (was 158968.5)
$ws.Cells.Item(84, 14).Value = -189981  # GSM!N84 (was -168952.5)
$ws.Cells.Item(86, 8).Value = 59999  # GSM!H86 (was 0)
$ws.Cells.Item(86, 10).Value = 59999  # GSM!J86 (was 0)
$ws.Cells.Item(86, 12).Value = 59999  # GSM!L86 (was 0)
$ws.Cells.Item(86, 14).Value = -62371  # GSM!N86 (was None)
$ws.Cells.Item(89, 8).Value = 59999  # GSM!H89 (was 0)
$ws.Cells.Item(89, 10).Value = 59999  # GSM!J89 (was 0)
$ws.Cells.Item(89, 12).Value = 179997  # GSM!L89 (was 0)
$ws.Cells.Item(89, 14).Value = -191853  # GSM!N89 (was None)

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 59900  # LTW!H6 (was 0)
$ws.Cells.Item(6, 10).Value = 59900  # LTW!J6 (was 0)
$ws.Cells.Item(6, 12).Value = 59900  # LTW!L6 (was 0)
$ws.Cells.Item(6, 14).Value = -60124  # LTW!N6 (was None)
$ws.Cells.Item(46, 8).Value = 1946.9412  # LTW!H46 (was 1832.1428)
$ws.Cells.Item(46, 9).Value = 2176.7693  # LTW!I46 (was 1977.2727)
$ws.Cells.Item(46, 10).Value = 1200  # LTW!J46 (was 1300)
$ws.Cells.Item(46, 11).Value = 2176.7693  # LTW!K46 (was 1977.2727)
$ws.Cells.Item(46, 12).Value = 1200  # LTW!L46 (was 1300)
$ws.Cells.Item(46, 13).Value = -1988.7693  # LTW!M46 (was -1789.2727)
$ws.Cells.Item(46, 14).Value = -1576  # LTW!N46 (was -1676)
$ws.Cells.Item(55, 8).Value = 397.125  # LTW!H55 (was 279.53845)
$ws.Cells.Item(55, 9).Value = 432  # LTW!I55 (was 277.18182)
$ws.Cells.Item(55, 11).Value = 432  # LTW!K55 (was 277.18182)
$ws.Cells.Item(55, 13).Value = -259  # LTW!M55 (was -104.18182)
$ws.Cells.Item(93, 8).Value = 0  # LTW!H93 (was 450)
$ws.Cells.Item(93, 10).Value = 0  # LTW!J93 (was 450)
$ws.Cells.Item(93, 12).Value = 0  # LTW!L93 (was 450)
$ws.Cells.Item(93, 14).Value = $null  # LTW!N93 clear (was -2946)
$ws.Cells.Item(94, 8).Value = 180000  # LTW!H94 (was 135000)
$ws.Cells.Item(94, 10).Value = 180000  # LTW!J94 (was 135000)
$ws.Cells.Item(94, 12).Value = 180000  # LTW!L94 (was 135000)
$ws.Cells.Item(94, 14).Value = -181352  # LTW!N94 (was -136352)
$ws.Cells.Item(98, 8).Value = 49998  # LTW!H98 (was 19000)
$ws.Cells.Item(98, 10).Value = 49998  # LTW!J98 (was 19000)
$ws.Cells.Item(98, 12).Value = 49998  # LTW!L98 (was 19000)
$ws.Cells.Item(98, 14).Value = -55988  # LTW!N98 (was -24990)

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 28999.666  # WVR!H2 (was 29000)
$ws.Cells.Item(2, 9).Value = 28999.5  # WVR!I2 (was 29000)
$ws.Cells.Item(2, 11).Value = 28999.5  # WVR!K2 (was 29000)
$ws.Cells.Item(2, 13).Value = -28887.5  # WVR!M2 (was -28888)
$ws.Cells.Item(74, 8).Value = 57801.5  # WVR!H74 (was 55041)
$ws.Cells.Item(74, 10).Value = 57801.5  # WVR!J74 (was 55041)
$ws.Cells.Item(74, 12).Value = 57801.5  # WVR!L74 (was 55041)
$ws.Cells.Item(74, 14).Value = -59673.5  # WVR!N74 (was -56913)
$ws.Cells.Item(77, 8).Value = 57801.5  # WVR!H77 (was 55041)
$ws.Cells.Item(77, 10).Value = 57801.5  # WVR!J77 (was 55041)
$ws.Cells.Item(77, 12).Value = 173404.5  # WVR!L77 (was 165123)
$ws.Cells.Item(77, 14).Value = -182764.5  # WVR!N77 (was -174483)
$ws.Cells.Item(81, 8).Value = 171165.83  # WVR!H81 (was 254250)
$ws.Cells.Item(81, 9).Value = 336666.66  # WVR!I81 (was 502500)
$ws.Cells.Item(81, 10).Value = 5665  # WVR!J81 (was 6000)
$ws.Cells.Item(81, 11).Value = 673333.32  # WVR!K81 (was 1005000)
$ws.Cells.Item(81, 12).Value = 11330  # WVR!L81 (was 12000)
$ws.Cells.Item(81, 13).Value = -672272.32  # WVR!M81 (was -1003939)
$ws.Cells.Item(81, 14).Value = -13452  # WVR!N81 (was -14122)
$ws.Cells.Item(84, 8).Value = 171165.83  # WVR!H84 (was 254250)
$ws.Cells.Item(84, 9).Value = 336666.66  # WVR!I84 (was 502500)
$ws.Cells.Item(84, 10).Value = 5665  # WVR!J84 (was 6000)
$ws.Cells.Item(84, 11).Value = 3366666.6  # WVR!K84 (was 5025000)
$ws.Cells.Item(84, 12).Value = 56650  # WVR!L84 (was 60000)
$ws.Cells.Item(84, 13).Value = -3361362.6  # WVR!M84 (was -5019696)
$ws.Cells.Item(84, 14).Value = -67258  # WVR!N84 (was -70608)
$ws.Cells.Item(136, 8).Value = 31251634  # WVR!H136 (was 31251640)
$ws.Cells.Item(136, 9).Value = 31251634  # WVR!I136 (was 31251640)
$ws.Cells.Item(136, 11).Value = 93754902  # WVR!K136 (was 93754920)
$ws.Cells.Item(136, 13).Value = -93752352  # WVR!M136 (was -93752370)
